$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A2").Value = '2025-10-08 12:37:38'
$ws.Range("B2").Value = '【Next.js × TypeScript × Tailwind】コンポーネント制作パートナー募集!'
$ws.Range("C2").Value = 'システム開発'
$ws.Range("D2").Value = '100,000 円 ~ 200,000 円 / 固定'
$ws.Range("E2").Value = '期限情報なし'
$ws.Range("F2").Value = 'https://www.lancers.jp/work/detail/5409073'
$ws.Range("G2").Value = 528
$ws.Range("H2").Value = '🔥AI,Next.js'

$ws.Range("A3").Value = '2025-10-08 12:37:38'
$ws.Range("B3").Value = '【急募】AI活用のRAGサービス開発をフルサポート!'
$ws.Range("C3").Value = 'システム開発'
$ws.Range("D3").Value = '1,000,000 円 ~ 3,000,000 円 / 固定'
$ws.Range("E3").Value = '期限情報なし'
$ws.Range("F3").Value = 'https://www.lancers.jp/work/detail/5409044'
$ws.Range("G3").Value = 375
$ws.Range("H3").Value = '🔥AI,Ai ◆開発'

$ws.Range("A4").Value = '2025-10-08 12:37:38'
$ws.Range("B4").Value = 'LLM+RAG活用の業務支援AI Agentサービスの構築を長期的に伴走できる開発パートナー募集'
$ws.Range("C4").Value = 'システム開発'
$ws.Range("D4").Value = '3,000,000 円 ~ 5,000,000 円 / 固定'
$ws.Range("E4").Value = '期限情報なし'
$ws.Range("F4").Value = 'https://www.lancers.jp/work/detail/5409015'
$ws.Range("G4").Value = 375
$ws.Range("H4").Value = '🔥AI,Ai ◆開発'

$ws.Range("A5").Value = '2025-10-08 12:37:38'
$ws.Range("B5").Value = '【急募】愛知県でのBtoB受発注システム開発者を募集!'
$ws.Range("C5").Value = 'システム開発'
$ws.Range("D5").Value = '50,000 円 ~ 100,000 円 / 固定'
$ws.Range("E5").Value = '期限情報なし'
$ws.Range("F5").Value = 'https://www.lancers.jp/work/detail/5408975'
$ws.Range("G5").Value = 118
$ws.Range("H5").Value = '◆開発,システム開発'

$ws.Range("A6").Value = '2025-10-08 12:37:38'
$ws.Range("B6").Value = '完全在宅GASエンジニア募集/課題テストからご依頼/時給1,163円~業務フロー効率化をお任せします'
$ws.Range("C6").Value = 'システム開発'
$ws.Range("D6").Value = '~ 5,000 円 / 固定'
$ws.Range("E6").Value = '期限情報なし'
$ws.Range("F6").Value = 'https://www.lancers.jp/work/detail/5409292'
$ws.Range("G6").Value = 70
$ws.Range("H6").Value = '◆効率化'

$ws.Range("A7").Value = '2025-10-08 12:37:38'
$ws.Range("B7").Value = '【急募】メールマガジンスタンド開発'
$ws.Range("C7").Value = 'システム開発'
$ws.Range("D7").Value = '50,000 円 ~ 100,000 円 / 固定'
$ws.Range("E7").Value = '期限情報なし'
$ws.Range("F7").Value = 'https://www.lancers.jp/work/detail/5409017'
$ws.Range("G7").Value = 68
$ws.Range("H7").Value = '◆開発'

$ws.Range("A8").Value = '2025-10-08 12:37:38'
$ws.Range("B8").Value = '【急募】wordpressでのECサイト更新作業をお手伝いください!'
$ws.Range("C8").Value = 'システム開発'
$ws.Range("D8").Value = '20,000 円 ~ 50,000 円 / 固定'
$ws.Range("E8").Value = '期限情報なし'
$ws.Range("F8").Value = 'https://www.lancers.jp/work/detail/5409217'
$ws.Range("G8").Value = 53
$ws.Range("H8").Value = '◇サイト ○WordPress'

$ws.Range("A9").Value = '2025-10-08 12:37:38'
$ws.Range("B9").Value = '【急募】Salesforce・MA・CRMコンサルタント経験者を探しています!'
$ws.Range("C9").Value = 'システム開発'
$ws.Range("D9").Value = '200,000 円 ~ 300,000 円 / 固定'
$ws.Range("E9").Value = '期限情報なし'
$ws.Range("F9").Value = 'https://www.lancers.jp/work/detail/5371747'
$ws.Range("G9").Value = 48
$ws.Range("H9").Value = '◆コンサル'

$ws.Range("A10").Value = '2025-10-08 12:37:38'
$ws.Range("B10").Value = 'AWSで運用中サイトの504及びクエリ数スパイクの原因調査'
$ws.Range("C10").Value = 'システム開発'
$ws.Range("D10").Value = '50,000 円 ~ 60,000 円 / 募集期間 3 日、取引期間 0 日'
$ws.Range("E10").Value = '期限情報なし'
$ws.Range("F10").Value = 'https://www.lancers.jp/work/detail/5408930'
$ws.Range("G10").Value = 33
$ws.Range("H10").Value = '◇サイト'

$ws.Range("A11").Value = '2025-10-08 12:37:38'
$ws.Range("B11").Value = '【時給2万円】iOSで他アプリ上に動画を重ねる仕組み(PiP等)を技術的に検証できるエンジニア募集'
$ws.Range("C11").Value = 'システム開発'
$ws.Range("D11").Value = '10,000 円 ~ 20,000 円 / 固定'
$ws.Range("E11").Value = '期限情報なし'
$ws.Range("F11").Value = 'https://www.lancers.jp/work/detail/5409283'
$ws.Range("G11").Value = 30
$ws.Range("H11").Value = '◇アプリ'

$ws.Range("A12").Value = '2025-10-08 12:37:38'
$ws.Range("B12").Value = '初回 運用中HPのドメイン分け'
$ws.Range("C12").Value = 'システム開発'
$ws.Range("D12").Value = '50,000 円 ~ 100,000 円 / 固定'
$ws.Range("E12").Value = '期限情報なし'
$ws.Range("F12").Value = 'https://www.lancers.jp/work/detail/5409114'
$ws.Range("G12").Value = 18
$ws.Range("H12").ClearContents()

$ws.Range("A13").Value = '2025-10-08 12:37:38'
$ws.Range("B13").Value = '【急募】Laravel/Filamentで構築したプログラム改修'
$ws.Range("C13").Value = 'システム開発'
$ws.Range("D13").Value = '50,000 円 ~ 100,000 円 / 固定'
$ws.Range("E13").Value = '期限情報なし'
$ws.Range("F13").Value = 'https://www.lancers.jp/work/detail/5409023'
$ws.Range("G13").Value = 18
$ws.Range("H13").ClearContents()

$ws.Range("A14").Value = '2025-10-08 12:37:38'
$ws.Range("B14").Value = '限定公開 限定公開の仕事'
$ws.Range("C14").Value = 'システム開発'
$ws.Range("D14").Value = '20,000 円 ~ 50,000 円 / 固定'
$ws.Range("E14").Value = '期限情報なし'
$ws.Range("F14").Value = 'https://www.lancers.jp/work/detail/5409366'
$ws.Range("G14").Value = 13
$ws.Range("H14").ClearContents()

$ws.Range("A15").Value = '2025-10-08 12:37:38'
$ws.Range("B15").Value = '【急募】家電商品の説明画像収集をお手伝いください!'
$ws.Range("C15").Value = 'システム開発'
$ws.Range("D15").Value = '20,000 円 ~ 50,000 円 / 固定'
$ws.Range("E15").Value = '期限情報なし'
$ws.Range("F15").Value = 'https://www.lancers.jp/work/detail/5409279'
$ws.Range("G15").Value = 13
$ws.Range("H15").ClearContents()

$ws.Range("A16").Value = '2025-10-08 12:37:38'
$ws.Range("B16").Value = 'bubbleでのサービス構築(difyとの連結)'
$ws.Range("C16").Value = 'システム開発'
$ws.Range("D16").Value = '20,000 円 ~ 50,000 円 / 固定'
$ws.Range("E16").Value = '期限情報なし'
$ws.Range("F16").Value = 'https://www.lancers.jp/work/detail/5408871'
$ws.Range("G16").Value = 13
$ws.Range("H16").ClearContents()

$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("F2"), 'https://www.lancers.jp/work/detail/5409073')
$ws.Hyperlinks.Add($ws.Range("F3"), 'https://www.lancers.jp/work/detail/5409044')
$ws.Hyperlinks.Add($ws.Range("F4"), 'https://www.lancers.jp/work/detail/5409015')
$ws.Hyperlinks.Add($ws.Range("F5"), 'https://www.lancers.jp/work/detail/5408975')
$ws.Hyperlinks.Add($ws.Range("F6"), 'https://www.lancers.jp/work/detail/5409292')
$ws.Hyperlinks.Add($ws.Range("F7"), 'https://www.lancers.jp/work/detail/5409017')
$ws.Hyperlinks.Add($ws.Range("F8"), 'https://www.lancers.jp/work/detail/5409217')
$ws.Hyperlinks.Add($ws.Range("F9"), 'https://www.lancers.jp/work/detail/5371747')
$ws.Hyperlinks.Add($ws.Range("F10"), 'https://www.lancers.jp/work/detail/5408930')
$ws.Hyperlinks.Add($ws.Range("F11"), 'https://www.lancers.jp/work/detail/5409283')
$ws.Hyperlinks.Add($ws.Range("F12"), 'https://www.lancers.jp/work/detail/5409114')
$ws.Hyperlinks.Add($ws.Range("F13"), 'https://www.lancers.jp/work/detail/5409023')
$ws.Hyperlinks.Add($ws.Range("F14"), 'https://www.lancers.jp/work/detail/5409366')
$ws.Hyperlinks.Add($ws.Range("F15"), 'https://www.lancers.jp/work/detail/5409279')
$ws.Hyperlinks.Add($ws.Range("F16"), 'https://www.lancers.jp/work/detail/5408871')

$ws.Range("F2:F16").Style = "Hyperlink"

